# Generate Report for Handoff
#
# A new handoff run produced a new source-file GUID (and new content
# hashes for the generated xliff files); refresh every cell/hyperlink
# in the report that is stamped with the old identifiers.

$wb = $excel.ActiveWorkbook

$oldGuid = "12979afa-ac46-43af-b4d3-87c3a31eac35"
$newGuid = "8923dd09-b6ad-4a18-87de-13e7141d813f"

$oldHash = "07317918ee3952b581a2df496c49c84235f3c5f0"
$newHash = "00e04a56a28372934ea03f0c9a48f89dc063a5ae"

$newOverviewDate = "2016-10-17 15:01:21"
$newZhHandoffDate = "2016-10-17 15:00:58"

# The hyperlink's underlying target URL is untouched by this commit -
# only the on-sheet display text changes - so reuse the existing address.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1060f6f3f0b72858453245517736228490b903e4/e2e/$oldGuid.md"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md")

$wsOverview.Range("G2").Value = $newOverviewDate

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, "$newGuid.md")

$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = $newZhHandoffDate

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, "$newGuid.md")

$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = $newOverviewDate
